$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new training session (CM "Introduction au langage C" on 2026-02-25 morning)
# was added to the planning, and an extra TP session (2026-03-06 afternoon,
# "Les boucles et les tableaux") was also added a bit further down. Net effect:
# one additional row in the schedule table (rows 2:13 -> rows 2:14).

# Insert a new blank row at row 7; this shifts the former rows 7:13 down to 8:14.
$ws.Rows("7").Insert()

# Give the new row 7 the same formatting as the (now shifted) row 8, which
# carries the "TP" style block used for this new entry.
$ws.Range("A8:D8").Copy()
$ws.Range("A7:D7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Now fill in the final values for every data row (A:D, rows 2:14) so that
# they match the updated planning.

# Row 2: 2026-02-25 (Wed) morning - CM - Introduction au langage C
$ws.Range("A2").Value = 46078
$ws.Range("B2").Value = 0.375
$ws.Range("C2").Value = "CM"
$ws.Range("D2").Value = "Introduction au langage C"

# Row 3: 2026-02-26 (Thu) morning - TP - Introduction au langage C
$ws.Range("A3").Value = 46079
$ws.Range("B3").Value = 0.375
$ws.Range("C3").Value = "TP"
$ws.Range("D3").Value = "Introduction au langage C"

# Row 4: 2026-02-26 (Thu) afternoon - CM - Fonctions, procedures et conditionnelle
$ws.Range("A4").Value = 46079
$ws.Range("B4").Value = 0.5625
$ws.Range("C4").Value = "CM"
$ws.Range("D4").Value = "Fonctions, procédures et conditionnelle"

# Row 5: 2026-03-04 (Wed) morning - TP - Fonctions, procedures et conditionnelle
$ws.Range("A5").Value = 46085
$ws.Range("B5").Value = 0.375
$ws.Range("C5").Value = "TP"
$ws.Range("D5").Value = "Fonctions, procédures et conditionnelle"

# Row 6: 2026-03-04 (Wed) afternoon - CM - Les boucles et les tableaux
$ws.Range("A6").Value = 46085
$ws.Range("B6").Value = 0.5625
$ws.Range("C6").Value = "CM"
$ws.Range("D6").Value = "Les boucles et les tableaux"

# Row 7: 2026-03-06 (Fri) morning - TP - Les boucles et les tableaux
$ws.Range("A7").Value = 46087
$ws.Range("B7").Value = 0.375
$ws.Range("C7").Value = "TP"
$ws.Range("D7").Value = "Les boucles et les tableaux"

# Row 8 (new): 2026-03-06 (Fri) afternoon - TP - Les boucles et les tableaux
$ws.Range("A8").Value = 46087
$ws.Range("B8").Value = 0.5625
$ws.Range("C8").Value = "TP"
$ws.Range("D8").Value = "Les boucles et les tableaux"

# Row 9: 2026-03-11 (Wed) morning - CM - Passage de parametres et fonctions
$ws.Range("A9").Value = 46092
$ws.Range("B9").Value = 0.375
$ws.Range("C9").Value = "CM"
$ws.Range("D9").Value = "Passage de paramètres et fonctions"

# Row 10: 2026-03-13 (Fri) morning - TP - Passage de parametres et fonctions
$ws.Range("A10").Value = 46094
$ws.Range("B10").Value = 0.375
$ws.Range("C10").Value = "TP"
$ws.Range("D10").Value = "Passage de paramètres et fonctions"

# Row 11: 2026-03-13 (Fri) afternoon - TP - Passage de parametres et fonctions
$ws.Range("A11").Value = 46094
$ws.Range("B11").Value = 0.5625
$ws.Range("C11").Value = "TP"
$ws.Range("D11").Value = "Passage de paramètres et fonctions"

# Row 12: 2026-03-18 (Wed) morning - TP - Mini Projet
$ws.Range("A12").Value = 46099
$ws.Range("B12").Value = 0.375
$ws.Range("C12").Value = "TP"
$ws.Range("D12").Value = "Mini Projet"

# Row 13: 2026-03-20 (Fri) morning - TP - Mini Projet
$ws.Range("A13").Value = 46101
$ws.Range("B13").Value = 0.375
$ws.Range("C13").Value = "TP"
$ws.Range("D13").Value = "Mini Projet"

# Row 14: 2026-03-25 (Wed) morning - DS - Examen
$ws.Range("A14").Value = 46106
$ws.Range("B14").Value = 0.375
$ws.Range("C14").Value = "DS"
$ws.Range("D14").Value = "Examen"

# Update the print area to cover the new row.
$ws.PageSetup.PrintArea = '$A$1:$D$14'

# Keep the active cell selection consistent with the edited workbook.
$ws.Range("E9").Select()
